$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-12-01 Monday"; new = "2025-12-02 Tuesday"},
    @{old = "763×7="; new = "202×6="},
    @{old = "544×4="; new = "219×6="},
    @{old = "536×7="; new = "877×9="},
    @{old = "369×9="; new = "530×9="},
    @{old = "254×5="; new = "798×5="},
    @{old = "325×2="; new = "722×2="},
    @{old = "253×4="; new = "587×2="},
    @{old = "230×3="; new = "423×6="},
    @{old = "223×8="; new = "224×2="},
    @{old = "418×9="; new = "619×5="},
    @{old = "207×4="; new = "931×5="},
    @{old = "988×5="; new = "544×9="},
    @{old = "333×2="; new = "481×9="},
    @{old = "728×5="; new = "625×2="},
    @{old = "843×2="; new = "180×2="},
    @{old = "463×2="; new = "458×6="},
    @{old = "746×8="; new = "530×9="},
    @{old = "721×8="; new = "309×9="},
    @{old = "582×8="; new = "454×5="},
    @{old = "345×9="; new = "130×3="},
    @{old = "508×5="; new = "494×6="},
    @{old = "616×9="; new = "839×3="},
    @{old = "434×7="; new = "354×3="},
    @{old = "308×9="; new = "398×3="},
    @{old = "750×4="; new = "673×9="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
